$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/eng-conversation-type"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Include from Engagement Conve" ---
$codes = $wb.Worksheets.Item("Include from Engagement Conve")
$codes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/eng-conversation-type"

$wb.Save()
